# Update countries & provincias Spain
# - Refresh "datos actualizados" timestamp
# - Reorder Mali / Maldivas rows (alphabetical: Maldivas before Mali) with refreshed
#   case counts for Maldivas (Mali's own figures are carried over unchanged)
# - Refresh case counts for Serbia, Kenia and Cabo Verde

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 15:34"

# Serbia (row 44): refreshed stats
$ws.Range("B44").Value = 9848
$ws.Range("C44").Value = 57
$ws.Range("D44").Value = 2160
$ws.Range("E44").Value = 7482
$ws.Range("F44").Value = 46
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 206

# Row 112 becomes "Maldivas" (alphabetically before "Mali") with refreshed stats
$ws.Range("A112").Value = "Maldivas"
$ws.Range("B112").Value = 642
$ws.Range("C112").Value = 25
$ws.Range("D112").Value = 20
$ws.Range("E112").Value = 620
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 2

# Row 113 becomes "Mali" carrying its previous (unchanged) stats
$ws.Range("A113").Value = "Mali"
$ws.Range("B113").Value = 631
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 261
$ws.Range("E113").Value = 338
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 32

# Kenia (row 116): refreshed stats
$ws.Range("B116").Value = 607
$ws.Range("C116").Value = 25
$ws.Range("D116").Value = 197
$ws.Range("E116").Value = 381
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = 29

# Cabo Verde (row 137): refreshed stats
$ws.Range("B137").Value = 218
$ws.Range("C137").Value = 27
$ws.Range("D137").Value = 38
$ws.Range("E137").Value = 178
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 2
